$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Helper pattern for Price (column D) cells: force the range to Text
# format before writing so Excel does not re-parse numeric-looking
# strings (e.g. "1.004", "0.00001085") into actual numbers, then put
# the style back to Normal so no stray formatting is left behind.
# ----------------------------------------------------------------------

# --- Refreshed Price / Volume(1h) figures (rows 2-43) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.217.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.793.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4489"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +15.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3699"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  +2.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07537"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.452"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.788.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001085"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "

$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("E21").Value = "  +3.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.349"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.191.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("E26").Value = "  +3.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.354"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.996.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.235"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09385"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.778"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2397"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02334"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06297"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.174"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6536"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.60%  "

$ws.Range("E42").Value = "  -1.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.198"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "

# --- Rows 44-47: coins reordered (Frax/EnergySwap swap, PancakeSwap/Decentraland swap) ---
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.62%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.821"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6055"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "

# --- Refreshed Price / Volume(1h) figures (rows 48-51) ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.023"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07117"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("E51").Value = "  +0.75%  "
